$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(16,2).Value2 = "CC"
$ws.Cells.Item(16,3).Value2 = "33102996"
$ws.Cells.Item(16,4).Value2 = "SHIRLIS ALVAREZ ESCALANTE"
$ws.Cells.Item(16,5).Value2 = "2201"
$ws.Cells.Item(16,6).Value2 = 31495
$ws.Cells.Item(16,7).Value2 = 908526
$ws.Cells.Item(17,2).Value2 = "CC"
$ws.Cells.Item(17,3).Value2 = "33102996"
$ws.Cells.Item(17,4).Value2 = "SHIRLIS ALVAREZ ESCALANTE"
$ws.Cells.Item(17,5).Value2 = "2112"
$ws.Cells.Item(17,6).Value2 = 36341
$ws.Cells.Item(17,7).Value2 = 908526
$ws.Cells.Item(18,2).Value2 = "CC"
$ws.Cells.Item(18,3).Value2 = "33102996"
$ws.Cells.Item(18,4).Value2 = "SHIRLIS ALVAREZ ESCALANTE"
$ws.Cells.Item(18,5).Value2 = "2111"
$ws.Cells.Item(18,6).Value2 = 36341
$ws.Cells.Item(18,7).Value2 = 908526
$ws.Cells.Item(19,2).Value2 = "CC"
$ws.Cells.Item(19,3).Value2 = "45460282"
$ws.Cells.Item(19,4).Value2 = "ALIEIDA DEL ROSARIO JULIO PINILLA"
$ws.Cells.Item(19,5).Value2 = "2201"
$ws.Cells.Item(19,6).Value2 = 24227
$ws.Cells.Item(19,7).Value2 = 908526
$ws.Cells.Item(20,2).Value2 = "CC"
$ws.Cells.Item(20,3).Value2 = "45460282"
$ws.Cells.Item(20,4).Value2 = "ALIEIDA DEL ROSARIO JULIO PINILLA"
$ws.Cells.Item(20,5).Value2 = "2112"
$ws.Cells.Item(20,6).Value2 = 36341
$ws.Cells.Item(20,7).Value2 = 908526
$ws.Cells.Item(21,2).Value2 = "CC"
$ws.Cells.Item(21,3).Value2 = "45460282"
$ws.Cells.Item(21,4).Value2 = "ALIEIDA DEL ROSARIO JULIO PINILLA"
$ws.Cells.Item(21,5).Value2 = "2111"
$ws.Cells.Item(21,6).Value2 = 36341
$ws.Cells.Item(21,7).Value2 = 908526
$ws.Cells.Item(22,2).Value2 = "CC"
$ws.Cells.Item(22,3).Value2 = "87880053"
$ws.Cells.Item(22,4).Value2 = "JUSTO TELLO"
$ws.Cells.Item(22,5).Value2 = "2201"
$ws.Cells.Item(22,6).Value2 = 104000
$ws.Cells.Item(22,7).Value2 = 3000000
$ws.Cells.Item(23,2).Value2 = "CC"
$ws.Cells.Item(23,3).Value2 = "87880053"
$ws.Cells.Item(23,4).Value2 = "JUSTO TELLO"
$ws.Cells.Item(23,5).Value2 = "2112"
$ws.Cells.Item(23,6).Value2 = 120000
$ws.Cells.Item(23,7).Value2 = 3000000
$ws.Cells.Item(24,2).Value2 = "CC"
$ws.Cells.Item(24,3).Value2 = "87880053"
$ws.Cells.Item(24,4).Value2 = "JUSTO TELLO"
$ws.Cells.Item(24,5).Value2 = "2111"
$ws.Cells.Item(24,6).Value2 = 120000
$ws.Cells.Item(24,7).Value2 = 3000000
$ws.Cells.Item(25,2).Value2 = "CC"
$ws.Cells.Item(25,3).Value2 = "1047385464"
$ws.Cells.Item(25,4).Value2 = "YOCELIN GOMEZ BARRIOS"
$ws.Cells.Item(25,5).Value2 = "2201"
$ws.Cells.Item(25,6).Value2 = 34666
$ws.Cells.Item(25,7).Value2 = 1300000
$ws.Cells.Item(26,2).Value2 = "CC"
$ws.Cells.Item(26,3).Value2 = "1047385464"
$ws.Cells.Item(26,4).Value2 = "YOCELIN GOMEZ BARRIOS"
$ws.Cells.Item(26,5).Value2 = "2112"
$ws.Cells.Item(26,6).Value2 = 52000
$ws.Cells.Item(26,7).Value2 = 1300000
$ws.Cells.Item(27,2).Value2 = "CC"
$ws.Cells.Item(27,3).Value2 = "1047385464"
$ws.Cells.Item(27,4).Value2 = "YOCELIN GOMEZ BARRIOS"
$ws.Cells.Item(27,5).Value2 = "2111"
$ws.Cells.Item(27,6).Value2 = 52000
$ws.Cells.Item(27,7).Value2 = 1300000
$ws.Cells.Item(28,2).Value2 = "CC"
$ws.Cells.Item(28,3).Value2 = "1002244348"
$ws.Cells.Item(28,4).Value2 = "FEDERICO DE JESUS DIAZ CASTRO"
$ws.Cells.Item(28,5).Value2 = "2201"
$ws.Cells.Item(28,6).Value2 = 34666
$ws.Cells.Item(28,7).Value2 = 1300000
$ws.Cells.Item(29,2).Value2 = "CC"
$ws.Cells.Item(29,3).Value2 = "1002244348"
$ws.Cells.Item(29,4).Value2 = "FEDERICO DE JESUS DIAZ CASTRO"
$ws.Cells.Item(29,5).Value2 = "2112"
$ws.Cells.Item(29,6).Value2 = 52000
$ws.Cells.Item(29,7).Value2 = 1300000
$ws.Cells.Item(30,2).Value2 = "CC"
$ws.Cells.Item(30,3).Value2 = "1002244348"
$ws.Cells.Item(30,4).Value2 = "FEDERICO DE JESUS DIAZ CASTRO"
$ws.Cells.Item(30,5).Value2 = "2111"
$ws.Cells.Item(30,6).Value2 = 52000
$ws.Cells.Item(30,7).Value2 = 1300000
$ws.Cells.Item(31,2).Value2 = "CC"
$ws.Cells.Item(31,3).Value2 = "30764183"
$ws.Cells.Item(31,4).Value2 = "GLORIA PATRICIA CASTRO CARRILLO"
$ws.Cells.Item(31,5).Value2 = "2201"
$ws.Cells.Item(31,6).Value2 = 34666
$ws.Cells.Item(31,7).Value2 = 1300000
$ws.Cells.Item(32,2).Value2 = "CC"
$ws.Cells.Item(32,3).Value2 = "30764183"
$ws.Cells.Item(32,4).Value2 = "GLORIA PATRICIA CASTRO CARRILLO"
$ws.Cells.Item(32,5).Value2 = "2112"
$ws.Cells.Item(32,6).Value2 = 52000
$ws.Cells.Item(32,7).Value2 = 1300000
$ws.Cells.Item(33,2).Value2 = "CC"
$ws.Cells.Item(33,3).Value2 = "30764183"
$ws.Cells.Item(33,4).Value2 = "GLORIA PATRICIA CASTRO CARRILLO"
$ws.Cells.Item(33,5).Value2 = "2111"
$ws.Cells.Item(33,6).Value2 = 52000
$ws.Cells.Item(33,7).Value2 = 1300000
